# Applies the scheduled runner update to all affected Leve profit sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 20 (HIKM updated)
$ws.Range("H20").Value = 21
$ws.Range("I20").Value = 21
$ws.Range("K20").Value = 21
$ws.Range("M20").Value = 209

# Row 33 (HIJKLM updated)
$ws.Range("H33").Value = 121.166664
$ws.Range("I33").Value = 121.166664
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 121.166664
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 107.833336
$ws.Range("N33").ClearContents()

# Row 35 (HIKM updated)
$ws.Range("H35").Value = 21
$ws.Range("I35").Value = 21
$ws.Range("K35").Value = 21
$ws.Range("M35").Value = 358

# Row 88 (HJLN updated)
$ws.Range("H88").Value = 19657894
$ws.Range("J88").Value = 56209.8
$ws.Range("L88").Value = 56209.8
$ws.Range("N88").Value = -57021.8

# Row 91 (HJLN updated)
$ws.Range("H91").Value = 19657894
$ws.Range("J91").Value = 56209.8
$ws.Range("L91").Value = 56209.8
$ws.Range("N91").Value = -59017.8

# Row 92 (HIKM updated)
$ws.Range("H92").Value = 31250594
$ws.Range("I92").Value = 491.81482
$ws.Range("K92").Value = 491.81482
$ws.Range("M92").Value = 756.1851799999999

# Row 113 (HJLN updated)
$ws.Range("H113").Value = 100708670
$ws.Range("J113").Value = 136383300
$ws.Range("L113").Value = 136383300
$ws.Range("N113").Value = -136389808

# Row 137 (HJLN updated)
$ws.Range("H137").Value = 6533.654
$ws.Range("J137").Value = 8425.866
$ws.Range("L137").Value = 25277.598
$ws.Range("N137").Value = -30377.598

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5 (HJLN updated)
$ws.Range("H5").Value = 30151
$ws.Range("J5").Value = 10002.5
$ws.Range("L5").Value = 10002.5
$ws.Range("N5").Value = -10226.5

# Row 88 (HIJKLMN updated)
$ws.Range("H88").Value = 2748.111
$ws.Range("I88").Value = 1006
$ws.Range("J88").Value = 2965.875
$ws.Range("K88").Value = 1006
$ws.Range("L88").Value = 2965.875
$ws.Range("M88").Value = -600
$ws.Range("N88").Value = -3777.875

# Row 91 (HIJKLMN updated)
$ws.Range("H91").Value = 2748.111
$ws.Range("I91").Value = 1006
$ws.Range("J91").Value = 2965.875
$ws.Range("K91").Value = 1006
$ws.Range("L91").Value = 2965.875
$ws.Range("M91").Value = 398
$ws.Range("N91").Value = -5773.875

# Row 138 (HJLN updated)
$ws.Range("H138").Value = 97232.2
$ws.Range("J138").Value = 97232.2
$ws.Range("L138").Value = 97232.2
$ws.Range("N138").Value = -107512.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4 (HJLN updated)
$ws.Range("H4").Value = 30151
$ws.Range("J4").Value = 10002.5
$ws.Range("L4").Value = 10002.5
$ws.Range("N4").Value = -10232.5

# Row 20 (HIJKLMN updated)
$ws.Range("H20").Value = 11118611
$ws.Range("I20").Value = 12825705
$ws.Range("J20").Value = 22500
$ws.Range("K20").Value = 12825705
$ws.Range("L20").Value = 22500
$ws.Range("M20").Value = -12825458
$ws.Range("N20").Value = -22994

# Row 105 (HIKM updated)
$ws.Range("H105").Value = 4884.5557
$ws.Range("I105").Value = 3975
$ws.Range("K105").Value = 3975
$ws.Range("M105").Value = -2228

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 2 (HJL updated)
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# Row 7 (HIKM updated)
$ws.Range("H7").Value = 239.125
$ws.Range("I7").Value = 157.9
$ws.Range("K7").Value = 157.9
$ws.Range("M7").Value = -44.90000000000001

# Row 16 (HIKM updated)
$ws.Range("H16").Value = 6697.4116
$ws.Range("I16").Value = 3663.6
$ws.Range("K16").Value = 3663.6
$ws.Range("M16").Value = -3376.6

# Row 58 (HIKM updated)
$ws.Range("H58").Value = 7301.147
$ws.Range("I58").Value = 2971.9092
$ws.Range("K58").Value = 2971.9092
$ws.Range("M58").Value = -2768.9092

# Row 63 (HJLN updated)
$ws.Range("H63").Value = 44663
$ws.Range("J63").Value = 44663
$ws.Range("L63").Value = 44663
$ws.Range("N63").Value = -46035

# Row 66 (HJLN updated)
$ws.Range("H66").Value = 44663
$ws.Range("J66").Value = 44663
$ws.Range("L66").Value = 133989
$ws.Range("N66").Value = -140853

# Row 69 (H updated)
$ws.Range("H69").Value = 32379.5

# Row 72 (H updated)
$ws.Range("H72").Value = 32379.5

# Row 99 (HJLN updated)
$ws.Range("H99").Value = 3977.739
$ws.Range("J99").Value = 5707.1816
$ws.Range("L99").Value = 5707.1816
$ws.Range("N99").Value = -8703.1816

# Row 103 (HIKM updated)
$ws.Range("H103").Value = 22435.5
$ws.Range("I103").Value = 11966.667
$ws.Range("K103").Value = 11966.667
$ws.Range("M103").Value = -10794.667

# Row 113 (HIKM updated)
$ws.Range("H113").Value = 6697.4116
$ws.Range("I113").Value = 3663.6
$ws.Range("K113").Value = 3663.6
$ws.Range("M113").Value = -1493.6

# Row 126 (HJLN updated)
$ws.Range("H126").Value = 3977.739
$ws.Range("J126").Value = 5707.1816
$ws.Range("L126").Value = 17121.5448
$ws.Range("N126").Value = -22061.5448

# Row 136 (HIKM updated)
$ws.Range("H136").Value = 7301.147
$ws.Range("I136").Value = 2971.9092
$ws.Range("K136").Value = 8915.7276
$ws.Range("M136").Value = -6365.7276

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2 (HIJKLMN updated)
$ws.Range("H2").Value = 96832.80499999999
$ws.Range("I2").Value = 18873.312
$ws.Range("J2").Value = 221568
$ws.Range("K2").Value = 113239.872
$ws.Range("L2").Value = 1329408
$ws.Range("M2").Value = -113126.872
$ws.Range("N2").Value = -1329634

# Row 4 (HIKM updated)
$ws.Range("H4").Value = 43070130
$ws.Range("I4").Value = 47832524
$ws.Range("K4").Value = 143497572
$ws.Range("M4").Value = -143497460

# Row 39 (HIJKLMN updated)
$ws.Range("H39").Value = 11488
$ws.Range("I39").Value = 11000
$ws.Range("J39").Value = 11610
$ws.Range("K39").Value = 33000
$ws.Range("L39").Value = 34830
$ws.Range("M39").Value = -32706
$ws.Range("N39").Value = -35418

# Row 98 (HIJKLMN updated)
$ws.Range("H98").Value = 2896.6
$ws.Range("I98").Value = 2003
$ws.Range("J98").Value = 3120
$ws.Range("K98").Value = 6009
$ws.Range("L98").Value = 9360
$ws.Range("M98").Value = -4511
$ws.Range("N98").Value = -12356

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70 (HIJKLMN updated)
$ws.Range("H70").Value = 5923.4
$ws.Range("I70").Value = 4624.364
$ws.Range("J70").Value = 7511.1113
$ws.Range("K70").Value = 4624.364
$ws.Range("L70").Value = 7511.1113
$ws.Range("M70").Value = -4354.364
$ws.Range("N70").Value = -8051.1113

# Row 73 (HIJKLMN updated)
$ws.Range("H73").Value = 5923.4
$ws.Range("I73").Value = 4624.364
$ws.Range("J73").Value = 7511.1113
$ws.Range("K73").Value = 4624.364
$ws.Range("L73").Value = 7511.1113
$ws.Range("M73").Value = -3688.364
$ws.Range("N73").Value = -9383.1113

# Row 80 (HJLN updated)
$ws.Range("H80").Value = 3104.8
$ws.Range("J80").Value = 3049.6667
$ws.Range("L80").Value = 3049.6667
$ws.Range("N80").Value = -5045.6667

# Row 83 (HJLN updated)
$ws.Range("H83").Value = 3104.8
$ws.Range("J83").Value = 3049.6667
$ws.Range("L83").Value = 15248.3335
$ws.Range("N83").Value = -25232.3335

# Row 97 (HIJKLMN updated)
$ws.Range("H97").Value = 1293.5834
$ws.Range("I97").Value = 657.5
$ws.Range("J97").Value = 2184.1
$ws.Range("K97").Value = 657.5
$ws.Range("L97").Value = 2184.1
$ws.Range("M97").Value = -161.5
$ws.Range("N97").Value = -3176.1

# Row 132 (HIJKLMN updated)
$ws.Range("H132").Value = 5361.3228
$ws.Range("I132").Value = 2870.1052
$ws.Range("J132").Value = 9305.75
$ws.Range("K132").Value = 8610.3156
$ws.Range("L132").Value = 27917.25
$ws.Range("M132").Value = -6080.3156
$ws.Range("N132").Value = -32977.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 20 (HJLN updated)
$ws.Range("H20").Value = 1070000.2
$ws.Range("J20").Value = 687500.25
$ws.Range("L20").Value = 687500.25
$ws.Range("N20").Value = -687952.25

# Row 122 (HIKM updated)
$ws.Range("H122").Value = 3207.3
$ws.Range("I122").Value = 2006.6428
$ws.Range("K122").Value = 6019.928400000001
$ws.Range("M122").Value = -3569.928400000001

# Row 132 (HIKM updated)
$ws.Range("H132").Value = 10422981
$ws.Range("I132").Value = 20836170
$ws.Range("K132").Value = 62508510
$ws.Range("M132").Value = -62505980

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100 (HIJKLMN updated)
$ws.Range("H100").Value = 552.96295
$ws.Range("I100").Value = 357.33334
$ws.Range("J100").Value = 944.2222
$ws.Range("K100").Value = 714.66668
$ws.Range("L100").Value = 1888.4444
$ws.Range("M100").Value = -173.66668
$ws.Range("N100").Value = -2970.4444

# Row 107 (HIJKLMN updated)
$ws.Range("H107").Value = 13889891
$ws.Range("I107").Value = 485.26666
$ws.Range("J107").Value = 37038900
$ws.Range("K107").Value = 1455.79998
$ws.Range("L107").Value = 111116700
$ws.Range("M107").Value = 464.20002
$ws.Range("N107").Value = -111120540

# Row 113 (HIJKLMN updated)
$ws.Range("H113").Value = 12052.833
$ws.Range("I113").Value = 22270.916
$ws.Range("J113").Value = 1834.75
$ws.Range("K113").Value = 66812.74800000001
$ws.Range("L113").Value = 5504.25
$ws.Range("M113").Value = -64642.74800000001
$ws.Range("N113").Value = -9844.25

# Row 135 (HJLN updated)
$ws.Range("H135").Value = 76499.75
$ws.Range("J135").Value = 76499.75
$ws.Range("L135").Value = 76499.75
$ws.Range("N135").Value = -86639.75

